$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells that lost their transcript text entirely
$ws.Range("C2").Value = ""
$ws.Range("C4").Value = ""
$ws.Range("C6").Value = ""

# Replace transcript text in the remaining data rows (rows 3,5,7-16)
$ws.Range("C3").Value = "全然誰から人口が増えるって事は家庭の子ども家庭の事を聞く人があってしまうので金銭的な都合で子供を作れたこともあると思う 製作 幼稚園から高校まで授業料無償化する"
$ws.Range("C5").Value = "インスト このテーマ めっちゃ 今の時代にいたの"
$ws.Range("C7").Value = "成増の高校とか後で教育費の無償化して回してもらう 子供を産んで待っていますけど逆にその子供がないとなんかにも税金がめちゃめちゃかかるみたいだから 一人一人の なんか人に対して なんかすごいで見ると終わるみたいな仕組みをもあるかと思います 世界にちんこを食い止めるから日本の人口 増えてるよ"
$ws.Range("C8").Value = "によって ちょっと変わってくるところ やけど あの難民と快眠とか日本とかに勧誘して 国籍を与えたら使用人口自体は増えることになるから日本から日本国内の人口だけで考えたらば 移民 難民とか受け入れ 国籍は連れて 3職業をつない 頭のやつをしてエッチしたら一応人口としたファイルがないで打つんが浮かんだかな 今コロナか 手間なくなっている人が増えてるから まぁ 結局的にワクチンを打って死亡していく人もなるべく 減少させていけば少しでも現場は 食い止められるんじゃないかと思います"
$ws.Range("C9").Value = "違かったコロナで死ぬより投資家 投資家の方が大問題としてもいいし そうしかをどうするか っていうの 考えると不眠にお見合い制度を設けるこれ これをやることで出会いがない人にも結婚したくない人以外は結婚できる確率は"
$ws.Range("C10").Value = "若者が住む独身で子供 結婚できないっていうのでその少子化が進むのでその明るくをもっと増やして住んでないとか お酒 小田君がいた出会いの場を設けてその環境はなくその環境とかその 言いやすくする環境 結婚しやすい環境なの作ればまあ 現象は抑えられるとそうでもなきゃ出来ないって 職場 そういう 稼げる場を出せばいいのかな"
$ws.Range("C11").Value = "今の意見 ある 仕方 最初の最初"
$ws.Range("C12").Value = "あと5分つあと5分 止めてもいいんだよいいけど"
$ws.Range("C13").Value = "求められる あと 幼稚園"
$ws.Range("C14").Value = "次は パチンコが増える 社会保険労働した 移民 難民とか仕事を与えてあげたら 下来てもいるしたチ**は安定した収入の安定収入を増やし 出会いの収入 面と出会いの場を設けて結婚をしやすくするという考え方 考えは 国が"
$ws.Range("C15").Value = "お見合いなくて結婚するきっかけ なども作ってそ歯科対策 お子様です 僕の考えはその制度的に例えばその教育費を無償化してまし F 2 でなくて 逆に独身の人とかにあの なんか 税金とかしてもらうようにしたほうがいい と僕 の考えは長期的に続いているので 僕的に摂取することを進めることです"
$ws.Range("C16").Value = "あと1分じゃ まとめると 次が政策を改めて何かしらの形で国民に援助金などを上げてもっと子供が暮らしやすい環境を作ったりと出会いの場を増やすことで子供を捨てる 引っ掛けなども 与えて 少子化を改善していくことで人口が 人口の減少を食い止めていけるのではないかという考えで大丈夫ですか"

# Rows 17-20 (15.wav..18.wav) are removed entirely from the sheet
$ws.Rows("17:20").Delete()
